{"js": "const pairs = [\n  [\"2025-09-09 Tuesday\", \"2025-09-10 Wednesday\"],\n  [\"188\u00f79=20, 8\", \"151\u00f74=37, 3\"],\n  [\"266\u00f74=66, 2\", \"287\u00f79=31, 8\"],\n  [\"760\u00f76=126, 4\", \"219\u00f72=109, 1\"],\n  [\"972\u00f76=162, 0\", \"248\u00f73=82, 2\"],\n  [\"282\u00f75=56, 2\", \"386\u00f72=193, 0\"],\n  [\"990\u00f79=110, 0\", \"469\u00f76=78, 1\"],\n  [\"957\u00f76=159, 3\", \"850\u00f77=121, 3\"],\n  [\"697\u00f74=174, 1\", \"701\u00f75=140, 1\"],\n  [\"900\u00f72=450, 0\", \"653\u00f75=130, 3\"],\n  [\"712\u00f72=356, 0\", \"646\u00f75=129, 1\"],\n  [\"833\u00f74=208, 1\", \"865\u00f72=432, 1\"],\n  [\"689\u00f73=229, 2\", \"774\u00f77=110, 4\"],\n  [\"445\u00f74=111, 1\", \"479\u00f75=95, 4\"],\n  [\"373\u00f77=53, 2\", \"188\u00f73=62, 2\"],\n  [\"719\u00f78=89, 7\", \"500\u00f77=71, 3\"],\n  [\"431\u00f74=107, 3\", \"320\u00f73=106, 2\"],\n  [\"144\u00f72=72, 0\", \"454\u00f77=64, 6\"],\n  [\"263\u00f75=52, 3\", \"286\u00f76=47, 4\"],\n  [\"695\u00f73=231, 2\", \"133\u00f78=16, 5\"],\n  [\"931\u00f72=465, 1\", \"579\u00f75=115, 4\"],\n  [\"848\u00f78=106, 0\", \"921\u00f76=153, 3\"],\n  [\"510\u00f74=127, 2\", \"257\u00f76=42, 5\"],\n  [\"958\u00f76=159, 4\", \"395\u00f77=56, 3\"],\n  [\"782\u00f79=86, 8\", \"141\u00f79=15, 6\"],\n  [\"313\u00f73=104, 1\", \"506\u00f78=63, 2\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $findText\n    $range.Find.Replacement.Text = $replaceText\n    $range.Find.Execute(\n        [ref]$findText,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$true,\n        [ref]1,\n        [ref]$false,\n        [ref]$replaceText,\n        [ref]2\n    ) | Out-Null\n}\n\nReplace-Text \"2025-09-09 Tuesday\" \"2025-09-10 Wednesday\"\nReplace-Text \"188\u00f79=20, 8\" \"151\u00f74=37, 3\"\nReplace-Text \"266\u00f74=66, 2\" \"287\u00f79=31, 8\"\nReplace-Text \"760\u00f76=126, 4\" \"219\u00f72=109, 1\"\nReplace-Text \"972\u00f76=162, 0\" \"248\u00f73=82, 2\"\nReplace-Text \"282\u00f75=56, 2\" \"386\u00f72=193, 0\"\nReplace-Text \"990\u00f79=110, 0\" \"469\u00f76=78, 1\"\nReplace-Text \"957\u00f76=159, 3\" \"850\u00f77=121, 3\"\nReplace-Text \"697\u00f74=174, 1\" \"701\u00f75=140, 1\"\nReplace-Text \"900\u00f72=450, 0\" \"653\u00f75=130, 3\"\nReplace-Text \"712\u00f72=356, 0\" \"646\u00f75=129, 1\"\nReplace-Text \"833\u00f74=208, 1\" \"865\u00f72=432, 1\"\nReplace-Text \"689\u00f73=229, 2\" \"774\u00f77=110, 4\"\nReplace-Text \"445\u00f74=111, 1\" \"479\u00f75=95, 4\"\nReplace-Text \"373\u00f77=53, 2\" \"188\u00f73=62, 2\"\nReplace-Text \"719\u00f78=89, 7\" \"500\u00f77=71, 3\"\nReplace-Text \"431\u00f74=107, 3\" \"320\u00f73=106, 2\"\nReplace-Text \"144\u00f72=72, 0\" \"454\u00f77=64, 6\"\nReplace-Text \"263\u00f75=52, 3\" \"286\u00f76=47, 4\"\nReplace-Text \"695\u00f73=231, 2\" \"133\u00f78=16, 5\"\nReplace-Text \"931\u00f72=465, 1\" \"579\u00f75=115, 4\"\nReplace-Text \"848\u00f78=106, 0\" \"921\u00f76=153, 3\"\nReplace-Text \"510\u00f74=127, 2\" \"257\u00f76=42, 5\"\nReplace-Text \"958\u00f76=159, 4\" \"395\u00f77=56, 3\"\nReplace-Text \"782\u00f79=86, 8\" \"141\u00f79=15, 6\"\nReplace-Text \"313\u00f73=104, 1\" \"506\u00f78=63, 2\"\n"}
